# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de)
# for the file "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md", mirroring the existing
# row for "831e25f3-a181-4e16-9403-096be0873547.md".

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c811ca4f26e5d1f7770f8b5db76227033e343556/e2e/"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"

$wsOverview.Range("C3").Value = ".md"

$wsOverview.Range("D3").Value = ""
$wsOverview.Range("D3").Font.Bold = $false

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsOverview.Range("G3").Value = "2016-09-05 00:45:35"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    ($repoBase + "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"),
    "",
    "",
    "e2e\cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"
) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"

$wsZhCn.Range("G3").Value = "cf967f90-bf43-4ca7-ba6c-79fc5358d722.2ab6e9d724a65e0cbe62cff4467ea2d68bb6874e.zh-cn.xlf"

$wsZhCn.Range("H3").Value = "2016-09-05 00:45:31"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("I3").Font.Bold = $false

$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("J3").Font.Bold = $false

$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("L3").Font.Bold = $false

$wsZhCn.Range("M3").Value = "'True"

$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("N3").Font.Bold = $false

$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Range("P3").Value = ""
$wsZhCn.Range("P3").Font.Bold = $false

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    ($repoBase + "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"),
    "",
    "",
    "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"
) | Out-Null
$wsZhCn.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"

$wsDeDe.Range("G3").Value = "cf967f90-bf43-4ca7-ba6c-79fc5358d722.2ab6e9d724a65e0cbe62cff4467ea2d68bb6874e.de-de.xlf"

$wsDeDe.Range("H3").Value = "2016-09-05 00:45:35"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("I3").Font.Bold = $false

$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("J3").Font.Bold = $false

$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("L3").Font.Bold = $false

$wsDeDe.Range("M3").Value = "'True"

$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("N3").Font.Bold = $false

$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Range("P3").Value = ""
$wsDeDe.Range("P3").Font.Bold = $false

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    ($repoBase + "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"),
    "",
    "",
    "cf967f90-bf43-4ca7-ba6c-79fc5358d722.md"
) | Out-Null
$wsDeDe.Range("A3").Style = "HyperLink"

Write-Host "Report row generated for handoff."
